$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: value correction ---
$ws1.Range("G2").Value = 0

# --- Sheet1: uncertainty columns (N) for rows 5,6,11,12 ---
$ws1.Range("N5").Value = 0.1
$ws1.Range("N6").Formula = "=SQRT((0.1)^2+(0.1)^2)"
$ws1.Range("N11").Value = 0.1
$ws1.Range("N12").Formula = "=SQRT((0.1)^2+(0.1)^2)"

# --- Sheet1: row 13 custom height, row 14 hidden ---
$ws1.Rows.Item(13).RowHeight = 14
$ws1.Rows.Item(14).Hidden = $true

# --- Sheet1: LINEST array formulas (rows 20-24) ---
$ws1.Range("J20:K24").FormulaArray = "=LINEST(B5:I5,B4:I4,TRUE,TRUE)"
$ws1.Range("M20:N24").FormulaArray = "=LINEST(C6:J6,B4:I4,TRUE,TRUE)"
$ws1.Range("P20:Q24").FormulaArray = "=LINEST(B8:E8,B4:E4,TRUE,TRUE)"

# --- Sheet1: LINEST array formulas (rows 27-31) ---
$ws1.Range("J27:K31").FormulaArray = "=LINEST(B11:K11,B4:K4,TRUE,TRUE)"
$ws1.Range("M27:N31").FormulaArray = "=LINEST(B12:K12,B4:K4,TRUE,TRUE)"
$ws1.Range("P27:Q31").FormulaArray = "=LINEST(B13:I13,B4:I4,TRUE,TRUE)"

# --- Sheet1: labels row 33 ("m" = slope label) ---
$ws1.Range("J33").Value = "m"
$ws1.Range("M33").Value = "m"
$ws1.Range("P33").Value = "m"

# --- Sheet1: label row 36 ("b" = intercept label) ---
$ws1.Range("J36").Value = "b"

# --- Sheet1: row 34 (c1 block, from slope m) ---
$ws1.Range("G34").Value = "c1"
$ws1.Range("J34").Formula = "=2*J20*M5"
$ws1.Range("K34").Formula = "=J34*SQRT((N11/M11)^2 + (J21/J20)^2)"
$ws1.Range("M34").Formula = "=M20*(2*M6)"
$ws1.Range("N34").Formula = "=M34*SQRT((N6/M6)^2 + (M21/M20)^2)"
$ws1.Range("P34").Formula = "=P20*M5"
$ws1.Range("Q34").Formula = "=P34*SQRT((N5/M5)^2 + (P21/P20)^2)"

# --- Sheet1: row 35 (c2 block, from slope m) ---
$ws1.Range("G35").Value = "c2"
$ws1.Range("J35").Formula = "=2*J27*M11"
$ws1.Range("K35").Formula = "=J35*SQRT((N11/M11)^2 + (J21/J20)^2)"
$ws1.Range("M35").Formula = "=M27*2*M12"
$ws1.Range("N35").Formula = "=M35*SQRT((N12/M12)^2 + (M28/M27)^2)"
$ws1.Range("P35").Formula = "=P27*M11"
$ws1.Range("Q35").Formula = "=P35*SQRT((N11/M11)^2 + (P28/P27)^2)"

# --- Sheet1: row 37 (c1 block, from intercept b) ---
$ws1.Range("G37").Value = "c1"
$ws1.Range("J37").Formula = "=-K20*4*M5"
$ws1.Range("K37").Formula = "=J37*SQRT((N5/M5)^2 + (K21/K20)^2)"
$ws1.Range("P37").Formula = "=-2*Q20*M5"
$ws1.Range("Q37").Formula = "=P37*SQRT((N5/M5)^2 + (Q21/Q20)^2)"

# --- Sheet1: row 38 (c2 block, from intercept b) ---
$ws1.Range("G38").Value = "c2"
$ws1.Range("J38").Formula = "=-K27*4*M11"
$ws1.Range("K38").Formula = "=J38*SQRT((N11/M11)^2 + (K28/K27)^2)"
$ws1.Range("P38").Formula = "=-2*Q27*M11"
$ws1.Range("Q38").Formula = "=P38*SQRT((N11/M11)^2 + (Q28/Q27)^2)"
$ws1.Range("J38").Interior.Color = 65535

# --- Sheet1: row 42 (weighted mean C1) ---
$ws1.Range("G42").Value = "C1 WM"
$ws1.Range("J42").Formula = "=J34/(K34^2)+M34/(N34^2)+P34/(Q34^2)+J37/(K37^2)+P37/(Q37^2)"
$ws1.Range("L42").Formula = "=1/(K34^2)+1/(N34^2)+1/(Q34^2)+1/(K37^2)+1/(Q37^2)"
$ws1.Range("M42").Formula = "=J42/L42"
$ws1.Range("N42").Formula = "=1/SQRT(L42)"

# --- Sheet1: row 43 (weighted mean C2) ---
$ws1.Range("G43").Value = "C2 WM"
$ws1.Range("J43").Formula = "=J35/(K35^2) + M35/(N35^2) + P35/(Q35^2) + J38/(K38^2) + P38/(Q38^2)"
$ws1.Range("L43").Formula = "=1/K35^2+1/N35^2+1/Q35^2+1/K38^2+1/Q38^2"
$ws1.Range("M43").Formula = "=J43/L43"
$ws1.Range("N43").Formula = "=1/SQRT(L43)"

# --- Sheet1: sheet view selection ---
$ws1.Range("F26").Select() | Out-Null

# --- Sheet2: selection ---
$ws2.Range("H17").Select() | Out-Null

Write-Output "done"
